# Adds a new "5. Evaluation" section (Heading1 + body paragraph) right
# before the trailing paragraph that carries the `_GoBack` bookmark,
# mirroring the target diff.

$d = $word.ActiveDocument

# --- Locate the paragraph that ends with "...distributions overlap...
# distribution." (the one immediately preceding the empty `_GoBack`
# bookmark paragraph at the tail of the document) via a text search, so
# we don't depend on a brittle hard-coded paragraph index. ---
$needle = "Boxplot analysis shows that score 5 is the main driver of significance"
$searchRange = $d.Content
$found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph for the Evaluation section."
}
$anchorPos = $searchRange.Start

$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $anchorPos -and $anchorPos -lt $pp.Range.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index."
}
$prevPara = $d.Paragraphs.Item($anchorIndex)

# Collapsed range positioned right before the anchor paragraph's own
# end-of-paragraph mark -- inserting WordOpenXML there cleanly creates
# brand-new sibling paragraphs after it (and before the following
# `_GoBack` paragraph) without disturbing either neighbour.
$insertAt = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End - 1)

$bodyText = "Collaboration among the team was superb and there were a well-defined line of communication and a proper allocation of tasks. We were able to formulate a policy-relevant research question that has practical implications on public health. The methodological rigor in the statistical analysis was proper in terms of the application of nonparametric statistics in non-normal distribution of data. The method of data visualization was successful in conveying complicated statistical relationships to different audiences. This extensive dataset allowed conducting systematic research in various health areas. Team meetings made sure that there were regular progress checks and that problem-solving within the team was done in good time."

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$headingPara = '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' + `
    '<w:bookmarkStart w:id="41" w:name="_Toc214911781"/>' + `
    '<w:bookmarkStart w:id="42" w:name="_Toc214915297"/>' + `
    '<w:bookmarkStart w:id="43" w:name="_Toc214915989"/>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>5. Evaluation</w:t></w:r>' + `
    '<w:bookmarkEnd w:id="41"/><w:bookmarkEnd w:id="42"/><w:bookmarkEnd w:id="43"/></w:p>'

$bodyPara = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>' + $bodyText + '</w:t></w:r></w:p>'

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $headingPara + $bodyPara + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertAt.InsertXML($xmlFragment)

Write-Output "Evaluation section inserted."
